# Refactor code to convert multiple tables:
# - rename Sheet1 to "dbo.Test"
# - add a new sheet "dbo.Test2" (a second, parallel table) right after it
# - make "dbo.Test2" the active sheet/tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "dbo.Test"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "dbo.Test2"

# Copy the whole table (values, types and styles) from the first sheet --
# the two tables share identical row data, only the header row differs.
$ws1.Range("A1:F5").Copy($ws2.Range("A1:F5"))

$ws2.Range("A1").Value = "Test2ID"
$ws2.Range("B1").Value = "Test2Content1"
$ws2.Range("C1").Value = "Test2Content2"
$ws2.Range("D1").Value = "Test2Content3"
$ws2.Range("E1").Value = "Test2Content4"
$ws2.Range("F1").Value = "Test2Content5"

# Match the column widths used by the original table.
$ws2.Columns.Item(1).ColumnWidth = 6.666666666666667
$ws2.Range("B1:F1").EntireColumn.ColumnWidth = 13.166666666666666

# Restore the selection on the first sheet (it is no longer the active tab).
$ws1.Range("A1:F5").Select()

# Select / activate the second sheet, matching the saved view state.
$ws2.Range("B2").Select()
$ws2.Activate()

Write-Host "done"
